$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in the source data (e.g. "42.338.61"),
# so a leading apostrophe forces Excel to keep them as text instead of parsing as numbers.

$ws.Range("D2").Value = "'42.338.61"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "'2.248.93"
$ws.Range("E3").Value = "  -3.78%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'236.23"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E6").Value = "  -4.60%  "
$ws.Range("D7").Value = "'69.73"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -6.62%  "
$ws.Range("D10").Value = "'0.0992"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "'58.81"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").Value = "'36.60"
$ws.Range("E12").Value = "  +13.47%  "
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "'6.77"
$ws.Range("E14").Value = "  -4.91%  "
$ws.Range("D15").Value = "'2.580.50"
$ws.Range("E15").Value = "  -3.89%  "
$ws.Range("E16").Value = "  -5.87%  "
$ws.Range("D17").Value = "'0.873"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").Value = "'2.247.85"
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").Value = "'42.258.76"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").Value = "'6.28"
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").Value = "'236.73"
$ws.Range("E23").Value = "  -5.60%  "
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  +7.60%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'3.66"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("E28").Value = "  +4.39%  "
$ws.Range("D29").Value = "'10.02"
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").Value = "'20.59"
$ws.Range("E31").Value = "  -6.85%  "
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("E33").Value = "  -4.47%  "
$ws.Range("D34").Value = "'0.0722"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "'4.72"
$ws.Range("E36").Value = "  -6.45%  "
$ws.Range("D37").Value = "'3.79"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").Value = "'23.13"
$ws.Range("E38").Value = "  +24.46%  "
$ws.Range("D39").Value = "'2.30"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").Value = "'0.0275"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("E41").Value = "  -6.16%  "
$ws.Range("D42").Value = "'65.99"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'9.41"
$ws.Range("E43").Value = "  +3.51%  "
$ws.Range("E44").Value = "  -15.22%  "
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("D46").Value = "'0.192"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'4.62"
$ws.Range("E47").Value = "  +13.59%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'10.25"
$ws.Range("E49").Value = "  +10.83%  "
$ws.Range("D51").Value = "'2.35"
$ws.Range("E51").Value = "  -2.04%  "
